$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 98000
$ws.Range("J95").Value = 98000
$ws.Range("L95").Value = 98000
$ws.Range("N95").Value = -103492

$ws.Range("H100").Value = 1829.7916
$ws.Range("I100").Value = 1791.591
$ws.Range("K100").Value = 1791.591
$ws.Range("M100").Value = -1250.591

$ws.Range("H139").Value = 73887.164
$ws.Range("J139").Value = 73887.164
$ws.Range("L139").Value = 73887.164
$ws.Range("N139").Value = -84167.164

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4774.9116
$ws.Range("I32").Value = 4032.4507
$ws.Range("J32").Value = 11364.25
$ws.Range("K32").Value = 4032.4507
$ws.Range("L32").Value = 11364.25
$ws.Range("M32").Value = -3745.4507
$ws.Range("N32").Value = -11938.25

$ws.Range("H45").Value = 1656
$ws.Range("I45").Value = 1552.1666
$ws.Range("J45").Value = 2071.3333
$ws.Range("K45").Value = 1552.1666
$ws.Range("L45").Value = 2071.3333
$ws.Range("M45").Value = -1175.1666
$ws.Range("N45").Value = -2825.3333

$ws.Range("H74").Value = 1984.6129
$ws.Range("I74").Value = 996.5
$ws.Range("J74").Value = 7122.8
$ws.Range("K74").Value = 996.5
$ws.Range("L74").Value = 7122.8
$ws.Range("M74").Value = -122.5
$ws.Range("N74").Value = -8870.799999999999

$ws.Range("H77").Value = 1984.6129
$ws.Range("I77").Value = 996.5
$ws.Range("J77").Value = 7122.8
$ws.Range("K77").Value = 4982.5
$ws.Range("L77").Value = 35614
$ws.Range("M77").Value = -614.5
$ws.Range("N77").Value = -44350

$ws.Range("H97").Value = 1396.0667
$ws.Range("I97").Value = 379.23077
$ws.Range("J97").Value = 8005.5
$ws.Range("K97").Value = 379.23077
$ws.Range("L97").Value = 8005.5
$ws.Range("M97").Value = 116.76923
$ws.Range("N97").Value = -8997.5

$ws.Range("H102").Value = 6581149.5
$ws.Range("I102").Value = 7355191
$ws.Range("J102").Value = 1800
$ws.Range("K102").Value = 7355191
$ws.Range("L102").Value = 1800
$ws.Range("M102").Value = -7353569
$ws.Range("N102").Value = -5044

$ws.Range("H110").Value = 884.2
$ws.Range("I110").Value = 731.9167
$ws.Range("K110").Value = 731.9167
$ws.Range("M110").Value = 1313.0833

$ws.Range("H132").Value = 4792.273
$ws.Range("I132").Value = 5847.1665
$ws.Range("J132").Value = 1979.2222
$ws.Range("K132").Value = 17541.4995
$ws.Range("L132").Value = 5937.6666
$ws.Range("M132").Value = -15011.4995
$ws.Range("N132").Value = -10997.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 825.10345
$ws.Range("I94").Value = 826.8461
$ws.Range("K94").Value = 826.8461
$ws.Range("M94").Value = -375.8461

$ws.Range("H99").Value = 1315.1904
$ws.Range("I99").Value = 1226.1875
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 1226.1875
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = 271.8125
$ws.Range("N99").Value = -4596

$ws.Range("H134").Value = 3405.6667
$ws.Range("I134").Value = 3797.1843
$ws.Range("J134").Value = 1917.9
$ws.Range("K134").Value = 11391.5529
$ws.Range("L134").Value = 5753.700000000001
$ws.Range("M134").Value = -8856.552899999999
$ws.Range("N134").Value = -10823.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18821.725
$ws.Range("I31").Value = 24845.809
$ws.Range("J31").Value = 3008.5
$ws.Range("K31").Value = 24845.809
$ws.Range("L31").Value = 3008.5
$ws.Range("M31").Value = -24550.809
$ws.Range("N31").Value = -3598.5

$ws.Range("H34").Value = 18821.725
$ws.Range("I34").Value = 24845.809
$ws.Range("J34").Value = 3008.5
$ws.Range("K34").Value = 24845.809
$ws.Range("L34").Value = 3008.5
$ws.Range("M34").Value = -24643.809
$ws.Range("N34").Value = -3412.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2998.2964
$ws.Range("I103").Value = 2145.4
$ws.Range("J103").Value = 3500
$ws.Range("K103").Value = 6436.200000000001
$ws.Range("L103").Value = 10500
$ws.Range("M103").Value = -5557.200000000001
$ws.Range("N103").Value = -12258

$ws.Range("H113").Value = 2975.8333
$ws.Range("I113").Value = 3411
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 10233
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = -8063
$ws.Range("N113").Value = -6740

$ws.Range("H134").Value = 4838
$ws.Range("I134").Value = 2176.5
$ws.Range("J134").Value = 7499.5
$ws.Range("K134").Value = 6529.5
$ws.Range("L134").Value = 22498.5
$ws.Range("M134").Value = -1459.5
$ws.Range("N134").Value = -32638.5

$ws.Range("H139").Value = 2682.7144
$ws.Range("I139").Value = 1044.75
$ws.Range("J139").Value = 4866.6665
$ws.Range("K139").Value = 3134.25
$ws.Range("L139").Value = 14599.9995
$ws.Range("M139").Value = 2005.75
$ws.Range("N139").Value = -24879.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 21748
$ws.Range("J95").Value = 21748
$ws.Range("L95").Value = 21748
$ws.Range("N95").Value = -27240

$ws.Range("H104").Value = 40671
$ws.Range("J104").Value = 40671
$ws.Range("L104").Value = 40671
$ws.Range("N104").Value = -47659

$ws.Range("H107").Value = 585.7368
$ws.Range("I107").Value = 236.55556
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 236.55556
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 1683.44444
$ws.Range("N107").Value = -4740

$ws.Range("H132").Value = 3704.3403
$ws.Range("I132").Value = 3804.5642
$ws.Range("J132").Value = 3215.75
$ws.Range("K132").Value = 11413.6926
$ws.Range("L132").Value = 9647.25
$ws.Range("M132").Value = -8883.692599999998
$ws.Range("N132").Value = -14707.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 11998
$ws.Range("J94").Value = 11998
$ws.Range("L94").Value = 11998
$ws.Range("N94").Value = -13350

$ws.Range("H136").Value = 3981.7847
$ws.Range("I136").Value = 2200.8386
$ws.Range("J136").Value = 5605.5884
$ws.Range("K136").Value = 6602.5158
$ws.Range("L136").Value = 16816.7652
$ws.Range("M136").Value = -4052.5158
$ws.Range("N136").Value = -21916.7652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 12122.429
$ws.Range("I43").Value = 4165.4
$ws.Range("J43").Value = 32015
$ws.Range("K43").Value = 4165.4
$ws.Range("L43").Value = 32015
$ws.Range("M43").Value = -4016.4
$ws.Range("N43").Value = -32313

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H136").Value = 7342.3125
$ws.Range("I136").Value = 8816.538
$ws.Range("J136").Value = 954
$ws.Range("K136").Value = 26449.614
$ws.Range("L136").Value = 2862
$ws.Range("M136").Value = -23899.614
$ws.Range("N136").Value = -7962
